# Apply F-column ("想去人数" / interested count) updates per the diff
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 234   # F3: 232 -> 234
$ws.Cells.Item(5, 6).Value = 9232   # F5: 9226 -> 9232
$ws.Cells.Item(6, 6).Value = 9232   # F6: 9226 -> 9232
$ws.Cells.Item(7, 6).Value = 561   # F7: 560 -> 561
$ws.Cells.Item(10, 6).Value = 235   # F10: 233 -> 235
$ws.Cells.Item(13, 6).Value = 147   # F13: 145 -> 147
$ws.Cells.Item(14, 6).Value = 151   # F14: 150 -> 151
$ws.Cells.Item(16, 6).Value = 11804   # F16: 11795 -> 11804
$ws.Cells.Item(17, 6).Value = 11804   # F17: 11795 -> 11804
$ws.Cells.Item(26, 6).Value = 21   # F26: 20 -> 21
$ws.Cells.Item(29, 6).Value = 2710   # F29: 2711 -> 2710
$ws.Cells.Item(32, 6).Value = 2093   # F32: 2092 -> 2093
$ws.Cells.Item(33, 6).Value = 59   # F33: 60 -> 59
$ws.Cells.Item(36, 6).Value = 973   # F36: 972 -> 973
$ws.Cells.Item(37, 6).Value = 4174   # F37: 4175 -> 4174
$ws.Cells.Item(39, 6).Value = 3596   # F39: 3595 -> 3596
$ws.Cells.Item(41, 6).Value = 2608   # F41: 2609 -> 2608
$ws.Cells.Item(42, 6).Value = 3048   # F42: 3049 -> 3048
$ws.Cells.Item(46, 6).Value = 403   # F46: 402 -> 403
$ws.Cells.Item(47, 6).Value = 467   # F47: 465 -> 467

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 13   # F6: 12 -> 13
$ws.Cells.Item(23, 6).Value = 32   # F23: 31 -> 32

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 6).Value = 13   # F5: 12 -> 13
$ws.Cells.Item(9, 6).Value = 234   # F9: 232 -> 234
$ws.Cells.Item(10, 6).Value = 9232   # F10: 9226 -> 9232
$ws.Cells.Item(11, 6).Value = 9232   # F11: 9226 -> 9232
$ws.Cells.Item(12, 6).Value = 561   # F12: 560 -> 561
$ws.Cells.Item(15, 6).Value = 235   # F15: 233 -> 235
$ws.Cells.Item(17, 6).Value = 147   # F17: 145 -> 147
$ws.Cells.Item(18, 6).Value = 151   # F18: 150 -> 151
$ws.Cells.Item(19, 6).Value = 11804   # F19: 11795 -> 11804
$ws.Cells.Item(20, 6).Value = 11804   # F20: 11795 -> 11804
$ws.Cells.Item(27, 6).Value = 21   # F27: 20 -> 21
$ws.Cells.Item(32, 6).Value = 2710   # F32: 2711 -> 2710
$ws.Cells.Item(35, 6).Value = 2093   # F35: 2092 -> 2093
$ws.Cells.Item(36, 6).Value = 59   # F36: 60 -> 59
$ws.Cells.Item(39, 6).Value = 973   # F39: 972 -> 973
$ws.Cells.Item(42, 6).Value = 3596   # F42: 3595 -> 3596
$ws.Cells.Item(43, 6).Value = 3048   # F43: 3049 -> 3048
$ws.Cells.Item(47, 6).Value = 403   # F47: 402 -> 403
$ws.Cells.Item(48, 6).Value = 32   # F48: 31 -> 32
$ws.Cells.Item(49, 6).Value = 467   # F49: 465 -> 467

$wb.Save()